$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new rows right after the existing "2040_TM152_FBP_Plus_20" row
# (row 85), pushing everything from the old row 86 onward down by two rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(86).Insert()
$ws.Rows.Item(86).Insert()

# ---------------------------------------------------------------------------
# Populate the two new rows (86 and 87) with the new FBP Plus runs for 2040.
# Columns: A=project, B=year, C=directory, D=run_set, E=category,
#          F=urbansim_path, G=urbansim_runid, H=status
# ---------------------------------------------------------------------------
$ws.Range("A86").Value = "RTP2021"
$ws.Range("B86").Value = 2040
$ws.Range("C86").Value = "2040_TM152_FBP_Plus_21"
$ws.Range("D86").Value = "FinalBlueprint"
$ws.Range("E86").Value = "Plus"
$ws.Range("F86").Value = "`"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION`""
$ws.Range("G86").Value = "run182"
$ws.Range("H86").Value = ""

$ws.Range("A87").Value = "RTP2021"
$ws.Range("B87").Value = 2040
$ws.Range("C87").Value = "2040_TM152_FBP_Plus_22"
$ws.Range("D87").Value = "FinalBlueprint"
$ws.Range("E87").Value = "Plus"
$ws.Range("F87").Value = "`"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION`""
$ws.Range("G87").Value = "run182"
$ws.Range("H87").Value = ""

# ---------------------------------------------------------------------------
# Fix up cell formatting/borders so that the "category / path / runid /
# status" mini-group (columns E:H) still brackets the correct set of rows:
#   - row 84 (NoProject_20) loses its bottom border on E:H (no longer last)
#   - row 85 (Plus_20) loses its bottom border on E:H (no longer last)
#   - row 86 (new Plus_21) has no border on E:H (middle row)
#   - row 87 (new Plus_22) gains the bottom border on E:H (now last row)
# Copy exact formats from existing donor cells that already carry the
# desired style so the result matches natively-authored formatting.
# ---------------------------------------------------------------------------
$ws.Range("E68:H68").Copy()
$ws.Range("E84").PasteSpecial(-4122) | Out-Null
$ws.Range("E85:H85").PasteSpecial(-4122) | Out-Null

$ws.Range("F15:H15").Copy()
$ws.Range("F84:H84").PasteSpecial(-4122) | Out-Null

$ws.Range("A15:H15").Copy()
$ws.Range("A86:H86").PasteSpecial(-4122) | Out-Null

$ws.Range("A15:D15").Copy()
$ws.Range("A87:D87").PasteSpecial(-4122) | Out-Null

$ws.Range("E19:H19").Copy()
$ws.Range("E87:H87").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Leave the selection on the newly added data, matching where the author
# was last working in the sheet.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A86").Select()
